$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-20 down to 8-21.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the weekly price entry (same market/category
# metadata as surrounding rows, new date + volume/price figures).
$ws.Cells.Item(7, 1).Value = 10
$ws.Cells.Item(7, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(7, 3).Value = "La Araucanía"
$ws.Cells.Item(7, 4).Value = 44720
$ws.Cells.Item(7, 5).Value = 9
$ws.Cells.Item(7, 6).Value = 100112042
$ws.Cells.Item(7, 7).Value = "Locoto"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 3600
$ws.Cells.Item(7, 12).Value = 3600
$ws.Cells.Item(7, 13).Value = 3600
$ws.Cells.Item(7, 14).Value = "$/kilo"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 3600
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date-time number format as the
# other rows in column D.
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
